$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 12.014481
$ws.Range("H2").Value = 36.043443
$ws.Range("I2").Value = 0.2338577750460931
$ws.Range("J2").Value = 0.2338577750460931
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 37.78741566666667
$ws.Range("N2").Value = 113.362247
$ws.Range("O2").Value = 0.3899441593213934
$ws.Range("P2").Value = 0.3899441593213934
$ws.Range("Q2").Value = 453.9961875662689
$ws.Range("R2").Value = 4085.96568809642
$ws.Range("S2").Value = 0.09119147349112032
$ws.Range("T2").Value = 0.0911914734911203

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 12.014481
$ws.Range("H3").Value = 36.043443
$ws.Range("I3").Value = 0.2338577750460931
$ws.Range("J3").Value = 0.2338577750460931
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 27.62253466666667
$ws.Range("N3").Value = 82.867604
$ws.Range("O3").Value = 0.2850484974663402
$ws.Range("P3").Value = 0.2850484974663403
$ws.Range("Q3").Value = 331.8704179245079
$ws.Range("R3").Value = 2986.833761320572
$ws.Range("S3").Value = 0.06666080739771024
$ws.Range("T3").Value = 0.06666080739771024

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 12.014481
$ws.Range("H4").Value = 36.043443
$ws.Range("I4").Value = 0.2338577750460931
$ws.Range("J4").Value = 0.2338577750460931
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 9.205342333333334
$ws.Range("N4").Value = 27.616027
$ws.Range("O4").Value = 0.09499378071000923
$ws.Range("P4").Value = 0.09499378071000925
$ws.Range("Q4").Value = 110.597410562329
$ws.Range("R4").Value = 995.376695060961
$ws.Range("S4").Value = 0.02221503420005924
$ws.Range("T4").Value = 0.02221503420005924

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 12.014481
$ws.Range("H5").Value = 36.043443
$ws.Range("I5").Value = 0.2338577750460931
$ws.Range("J5").Value = 0.2338577750460931
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 17.41884833333333
$ws.Range("N5").Value = 52.256545
$ws.Range("O5").Value = 0.1797523871334833
$ws.Range("P5").Value = 0.1797523871334834
$ws.Range("Q5").Value = 209.278422342715
$ws.Range("R5").Value = 1883.505801084435
$ws.Range("S5").Value = 0.0420364933142604
$ws.Range("T5").Value = 0.0420364933142604

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 12.014481
$ws.Range("H6").Value = 36.043443
$ws.Range("I6").Value = 0.2338577750460931
$ws.Range("J6").Value = 0.2338577750460931
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 4.870543333333333
$ws.Range("N6").Value = 14.61163
$ws.Range("O6").Value = 0.05026117536877379
$ws.Range("P6").Value = 0.0502611753687738
$ws.Range("Q6").Value = 58.51705033800999
$ws.Range("R6").Value = 526.65345304209
$ws.Range("S6").Value = 0.01175396664294294
$ws.Range("T6").Value = 0.01175396664294294

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 11.16283
$ws.Range("H7").Value = 33.48849
$ws.Range("I7").Value = 0.217280678792349
$ws.Range("J7").Value = 0.217280678792349
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 37.78741566666667
$ws.Range("N7").Value = 113.362247
$ws.Range("O7").Value = 0.3899441593213934
$ws.Range("P7").Value = 0.3899441593213934
$ws.Range("Q7").Value = 421.8144972263366
$ws.Range("R7").Value = 3796.33047503703
$ws.Range("S7").Value = 0.08472733162846423
$ws.Range("T7").Value = 0.08472733162846424

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 11.16283
$ws.Range("H8").Value = 33.48849
$ws.Range("I8").Value = 0.217280678792349
$ws.Range("J8").Value = 0.217280678792349
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 27.62253466666667
$ws.Range("N8").Value = 82.867604
$ws.Range("O8").Value = 0.2850484974663402
$ws.Range("P8").Value = 0.2850484974663403
$ws.Range("Q8").Value = 308.3456586531067
$ws.Range("R8").Value = 2775.11092787796
$ws.Range("S8").Value = 0.06193553101822557
$ws.Range("T8").Value = 0.06193553101822558

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 11.16283
$ws.Range("H9").Value = 33.48849
$ws.Range("I9").Value = 0.217280678792349
$ws.Range("J9").Value = 0.217280678792349
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 9.205342333333334
$ws.Range("N9").Value = 27.616027
$ws.Range("O9").Value = 0.09499378071000923
$ws.Range("P9").Value = 0.09499378071000925
$ws.Range("Q9").Value = 102.7576715588033
$ws.Range("R9").Value = 924.8190440292301
$ws.Range("S9").Value = 0.02064031315372235
$ws.Range("T9").Value = 0.02064031315372235

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 11.16283
$ws.Range("H10").Value = 33.48849
$ws.Range("I10").Value = 0.217280678792349
$ws.Range("J10").Value = 0.217280678792349
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 17.41884833333333
$ws.Range("N10").Value = 52.256545
$ws.Range("O10").Value = 0.1797523871334833
$ws.Range("P10").Value = 0.1797523871334834
$ws.Range("Q10").Value = 194.4436427407833
$ws.Range("R10").Value = 1749.99278466705
$ws.Range("S10").Value = 0.03905672069090836
$ws.Range("T10").Value = 0.03905672069090836

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 11.16283
$ws.Range("H11").Value = 33.48849
$ws.Range("I11").Value = 0.217280678792349
$ws.Range("J11").Value = 0.217280678792349
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 4.870543333333333
$ws.Range("N11").Value = 14.61163
$ws.Range("O11").Value = 0.05026117536877379
$ws.Range("P11").Value = 0.0502611753687738
$ws.Range("Q11").Value = 54.36904723763332
$ws.Range("R11").Value = 489.3214251387
$ws.Range("S11").Value = 0.01092078230102846
$ws.Range("T11").Value = 0.01092078230102846

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 13.04190333333333
$ws.Range("H12").Value = 39.12571
$ws.Range("I12").Value = 0.2538562003551846
$ws.Range("J12").Value = 0.2538562003551846
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 37.78741566666667
$ws.Range("N12").Value = 113.362247
$ws.Range("O12").Value = 0.3899441593213934
$ws.Range("P12").Value = 0.3899441593213934
$ws.Range("Q12").Value = 492.8198223411522
$ws.Range("R12").Value = 4435.378401070369
$ws.Range("S12").Value = 0.09898974263602567
$ws.Range("T12").Value = 0.09898974263602567

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 13.04190333333333
$ws.Range("H13").Value = 39.12571
$ws.Range("I13").Value = 0.2538562003551846
$ws.Range("J13").Value = 0.2538562003551846
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 27.62253466666667
$ws.Range("N13").Value = 82.867604
$ws.Range("O13").Value = 0.2850484974663402
$ws.Range("P13").Value = 0.2850484974663403
$ws.Range("Q13").Value = 360.2504269443155
$ws.Range("R13").Value = 3242.25384249884
$ws.Range("S13").Value = 0.0723613284837596
$ws.Range("T13").Value = 0.07236132848375959

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 13.04190333333333
$ws.Range("H14").Value = 39.12571
$ws.Range("I14").Value = 0.2538562003551846
$ws.Range("J14").Value = 0.2538562003551846
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 9.205342333333334
$ws.Range("N14").Value = 27.616027
$ws.Range("O14").Value = 0.09499378071000923
$ws.Range("P14").Value = 0.09499378071000925
$ws.Range("Q14").Value = 120.0551848615744
$ws.Range("R14").Value = 1080.49666375417
$ws.Range("S14").Value = 0.02411476022841658
$ws.Range("T14").Value = 0.02411476022841658

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 13.04190333333333
$ws.Range("H15").Value = 39.12571
$ws.Range("I15").Value = 0.2538562003551846
$ws.Range("J15").Value = 0.2538562003551846
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 17.41884833333333
$ws.Range("N15").Value = 52.256545
$ws.Range("O15").Value = 0.1797523871334833
$ws.Range("P15").Value = 0.1797523871334834
$ws.Range("Q15").Value = 227.1749361413277
$ws.Range("R15").Value = 2044.57442527195
$ws.Range("S15").Value = 0.04563125800248026
$ws.Range("T15").Value = 0.04563125800248026

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 13.04190333333333
$ws.Range("H16").Value = 39.12571
$ws.Range("I16").Value = 0.2538562003551846
$ws.Range("J16").Value = 0.2538562003551846
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 4.870543333333333
$ws.Range("N16").Value = 14.61163
$ws.Range("O16").Value = 0.05026117536877379
$ws.Range("P16").Value = 0.0502611753687738
$ws.Range("Q16").Value = 63.52115533414443
$ws.Range("R16").Value = 571.6903980073
$ws.Range("S16").Value = 0.01275911100450251
$ws.Range("T16").Value = 0.01275911100450251

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 11.62036266666667
$ws.Range("H17").Value = 34.861088
$ws.Range("I17").Value = 0.226186396104447
$ws.Range("J17").Value = 0.2261863961044469
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 37.78741566666667
$ws.Range("N17").Value = 113.362247
$ws.Range("O17").Value = 0.3899441593213934
$ws.Range("P17").Value = 0.3899441593213934
$ws.Range("Q17").Value = 439.1034742827484
$ws.Range("R17").Value = 3951.931268544735
$ws.Range("S17").Value = 0.08820006407888427
$ws.Range("T17").Value = 0.08820006407888425

$ws.Range("E18").Value = 3
$ws.Range("G18").Value = 11.62036266666667
$ws.Range("H18").Value = 34.861088
$ws.Range("I18").Value = 0.226186396104447
$ws.Range("J18").Value = 0.2261863961044469
$ws.Range("K18").Value = 3
$ws.Range("M18").Value = 27.62253466666667
$ws.Range("N18").Value = 82.867604
$ws.Range("O18").Value = 0.2850484974663402
$ws.Range("P18").Value = 0.2850484974663403
$ws.Range("Q18").Value = 320.9838705992391
$ws.Range("R18").Value = 2888.854835393151
$ws.Range("S18").Value = 0.06447409235689908
$ws.Range("T18").Value = 0.06447409235689908

$ws.Range("E19").Value = 3
$ws.Range("G19").Value = 11.62036266666667
$ws.Range("H19").Value = 34.861088
$ws.Range("I19").Value = 0.226186396104447
$ws.Range("J19").Value = 0.2261863961044469
$ws.Range("K19").Value = 3
$ws.Range("M19").Value = 9.205342333333334
$ws.Range("N19").Value = 27.616027
$ws.Range("O19").Value = 0.09499378071000923
$ws.Range("P19").Value = 0.09499378071000925
$ws.Range("Q19").Value = 106.9694163841529
$ws.Range("R19").Value = 962.724747457376
$ws.Range("S19").Value = 0.02148630091113312
$ws.Range("T19").Value = 0.02148630091113312

$ws.Range("E20").Value = 3
$ws.Range("G20").Value = 11.62036266666667
$ws.Range("H20").Value = 34.861088
$ws.Range("I20").Value = 0.226186396104447
$ws.Range("J20").Value = 0.2261863961044469
$ws.Range("K20").Value = 3
$ws.Range("M20").Value = 17.41884833333333
$ws.Range("N20").Value = 52.256545
$ws.Range("O20").Value = 0.1797523871334833
$ws.Range("P20").Value = 0.1797523871334834
$ws.Range("Q20").Value = 202.4133348689955
$ws.Range("R20").Value = 1821.72001382096
$ws.Range("S20").Value = 0.04065754463689396
$ws.Range("T20").Value = 0.04065754463689396

$ws.Range("E21").Value = 3
$ws.Range("G21").Value = 11.62036266666667
$ws.Range("H21").Value = 34.861088
$ws.Range("I21").Value = 0.226186396104447
$ws.Range("J21").Value = 0.2261863961044469
$ws.Range("K21").Value = 3
$ws.Range("M21").Value = 4.870543333333333
$ws.Range("N21").Value = 14.61163
$ws.Range("O21").Value = 0.05026117536877379
$ws.Range("P21").Value = 0.0502611753687738
$ws.Range("Q21").Value = 56.59747991704888
$ws.Range("R21").Value = 509.3773192534399
$ws.Range("S21").Value = 0.01136839412063654
$ws.Range("T21").Value = 0.01136839412063654

$ws.Range("E22").Value = 3
$ws.Range("G22").Value = 3.535584666666667
$ws.Range("H22").Value = 10.606754
$ws.Range("I22").Value = 0.06881894970192634
$ws.Range("J22").Value = 0.06881894970192633
$ws.Range("K22").Value = 3
$ws.Range("M22").Value = 37.78741566666667
$ws.Range("N22").Value = 113.362247
$ws.Range("O22").Value = 0.3899441593213934
$ws.Range("P22").Value = 0.3899441593213934
$ws.Range("Q22").Value = 133.6006074240265
$ws.Range("R22").Value = 1202.405466816238
$ws.Range("S22").Value = 0.02683554748689892
$ws.Range("T22").Value = 0.02683554748689892

$ws.Range("E23").Value = 3
$ws.Range("G23").Value = 3.535584666666667
$ws.Range("H23").Value = 10.606754
$ws.Range("I23").Value = 0.06881894970192634
$ws.Range("J23").Value = 0.06881894970192633
$ws.Range("K23").Value = 3
$ws.Range("M23").Value = 27.62253466666667
$ws.Range("N23").Value = 82.867604
$ws.Range("O23").Value = 0.2850484974663402
$ws.Range("P23").Value = 0.2850484974663403
$ws.Range("Q23").Value = 97.66181002193511
$ws.Range("R23").Value = 878.956290197416
$ws.Range("S23").Value = 0.01961673820974574
$ws.Range("T23").Value = 0.01961673820974574

$ws.Range("E24").Value = 3
$ws.Range("G24").Value = 3.535584666666667
$ws.Range("H24").Value = 10.606754
$ws.Range("I24").Value = 0.06881894970192634
$ws.Range("J24").Value = 0.06881894970192633
$ws.Range("K24").Value = 3
$ws.Range("M24").Value = 9.205342333333334
$ws.Range("N24").Value = 27.616027
$ws.Range("O24").Value = 0.09499378071000923
$ws.Range("P24").Value = 0.09499378071000925
$ws.Range("Q24").Value = 32.54626720515089
$ws.Range("R24").Value = 292.916404846358
$ws.Range("S24").Value = 0.006537372216677946
$ws.Range("T24").Value = 0.006537372216677946

$ws.Range("E25").Value = 3
$ws.Range("G25").Value = 3.535584666666667
$ws.Range("H25").Value = 10.606754
$ws.Range("I25").Value = 0.06881894970192634
$ws.Range("J25").Value = 0.06881894970192633
$ws.Range("K25").Value = 3
$ws.Range("M25").Value = 17.41884833333333
$ws.Range("N25").Value = 52.256545
$ws.Range("O25").Value = 0.1797523871334833
$ws.Range("P25").Value = 0.1797523871334834
$ws.Range("Q25").Value = 61.58581307832556
$ws.Range("R25").Value = 554.2723177049301
$ws.Range("S25").Value = 0.01237037048894038
$ws.Range("T25").Value = 0.01237037048894038

$ws.Range("E26").Value = 3
$ws.Range("G26").Value = 3.535584666666667
$ws.Range("H26").Value = 10.606754
$ws.Range("I26").Value = 0.06881894970192634
$ws.Range("J26").Value = 0.06881894970192633
$ws.Range("K26").Value = 3
$ws.Range("M26").Value = 4.870543333333333
$ws.Range("N26").Value = 14.61163
$ws.Range("O26").Value = 0.05026117536877379
$ws.Range("P26").Value = 0.0502611753687738
$ws.Range("Q26").Value = 17.22021832766889
$ws.Range("R26").Value = 154.98196494902
$ws.Range("S26").Value = 0.003458921299663343
$ws.Range("T26").Value = 0.003458921299663343
